$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new day's data to row 18
$ws.Range("A18").Value = 45960
$ws.Range("B18").Value = 5602
$ws.Range("C18").Value = 4316
$ws.Range("D18").Value = 3961
$ws.Range("E18").Value = 265
$ws.Range("F18").Value = 55
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 1

# Move the current selection down to the newly filled row
$ws.Range("A18:I18").Select()
